# This edit re-orders the species-observation records that live in rows 3-7
# of the sheet (a 5-cycle permutation of the rows), while columns C, T, U, V,
# W, Z, AB, AD, AE, AG, AT, AW, AY stay constant across all of these rows.
#
# Mapping of new row -> original row that its data came from:
#   row 3 <- old row 4
#   row 4 <- old row 5
#   row 5 <- old row 7
#   row 6 <- old row 3
#   row 7 <- old row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes one text cell while keeping it as literal text instead
# of letting Excel auto-convert date-looking strings (e.g. "2020-10-04")
# into a date serial number.
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

function Set-RowData($ws, $row, $A, $B, $D, $E, $F, $G, $H, $P, $Q, $R, $Y, $AA, $AX) {
    $ws.Range("A$row").Value = $A
    $ws.Range("B$row").Value = $B
    $ws.Range("D$row").Value = $D
    $ws.Range("E$row").Value = $E
    $ws.Range("F$row").Value = $F
    $ws.Range("G$row").Value = $G
    $ws.Range("H$row").Value = $H
    $ws.Range("P$row").Value = $P
    $ws.Range("Q$row").Value = $Q
    $ws.Range("R$row").Value = $R
    Set-TextCell $ws "Y$row" $Y
    Set-TextCell $ws "AA$row" $AA
    $ws.Range("AX$row").Value = $AX
}

# New row 3 (was old row 4)
Set-RowData $ws 3 88608218 78603 "LC" 6464 "Luddlav" "Nephroma resupinatum" "(L.) Ach." "Kvilåsen, Jmt" 440619.8703747808 7163001.094314476 "2020-10-04" "2020-10-04" "Kristina Bäck"

# New row 4 (was old row 5)
Set-RowData $ws 4 88608217 78596 "LC" 6462 "Stuplav" "Nephroma bellum" "(Spreng.) Tuck." "Kvilåsen, Jmt" 440619.8703747808 7163001.094314476 "2020-10-04" "2020-10-04" "Kristina Bäck"

# New row 5 (was old row 7)
Set-RowData $ws 5 88608204 78570 "NT" 2081 "Skrovellav" "Lobaria scrobiculata" "(Scop.) DC." "Kvilåsen, Jmt" 440666.2172450395 7163281.207762145 "2020-10-04" "2020-10-04" "Kristina Bäck"

# New row 6 (was old row 3)
Set-RowData $ws 6 88608299 78596 "LC" 6462 "Stuplav" "Nephroma bellum" "(Spreng.) Tuck." "Björkvattsruet, Jmt" 440455.1034245967 7162211.825440676 "2020-10-03" "2020-10-03" "Via Maria Danvind"

# New row 7 (was old row 6)
Set-RowData $ws 7 88608297 77506 "NT" 6425 "Garnlav" "Alectoria sarmentosa" "(Ach.) Ach." "Björkvattsruet, Jmt" 440455.1034245967 7162211.825440676 "2020-10-03" "2020-10-03" "Via Maria Danvind"
